# Natmi LR-pair output following Dr Hou's advice: the mean/total expression
# columns were recomputed per sending/target cluster pair, and the cluster
# set used for "Target cluster" now also covers "ECs" as a Sending cluster,
# producing the full 3x3 cross (ECs / FAPs / sCs) x (ECs / FAPs / sCs) for
# the Angpt1->Itga5 pair (9 data rows instead of 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol |
#          D Target cluster  | E..T the numeric metric columns
$rows = @(
    ,@("ECs",  "Angpt1", "Itga5", "ECs",  2, 0.6666666666666666, 0.1659303333333333, 0.497791,           0.01127931880021135, 0.01127931880021135, 3, 1, 34.07074633333333, 102.212239,        0.5171464495142372, 0.5171464495142373, 5.653370296005442,  50.88033266404899,  0.005833059670468486, 0.005833059670468487)
    ,@("ECs",  "Angpt1", "Itga5", "FAPs", 2, 0.6666666666666666, 0.1659303333333333, 0.497791,           0.01127931880021135, 0.01127931880021135, 3, 1, 27.685497,          83.056491,         0.420227262899125,  0.4202272628991251, 4.593863745708999,  41.34477371138099,  0.00473987726677946,  0.00473987726677946)
    ,@("ECs",  "Angpt1", "Itga5", "sCs",  2, 0.6666666666666666, 0.1659303333333333, 0.497791,           0.01127931880021135, 0.01127931880021135, 3, 1, 4.125957666666666, 12.377873,         0.06262628758663766,0.06262628758663766,0.6846215309492221, 6.161593778543,     0.0007063818629634049,0.0007063818629634049)
    ,@("FAPs", "Angpt1", "Itga5", "ECs",  3, 1,                  12.95459633333333, 38.863789,           0.8806046431437031,  0.8806046431437031,  3, 1, 34.07074633333333, 102.212239,        0.5171464495142372, 0.5171464495142373, 441.37276552373,    3972.35488971357,   0.4554015646275179,   0.455401564627518)
    ,@("FAPs", "Angpt1", "Itga5", "FAPs", 3, 1,                  12.95459633333333, 38.863789,           0.8806046431437031,  0.8806046431437031,  3, 1, 27.685497,          83.056491,         0.420227262899125,  0.4202272628991251, 358.6544379227109,  3227.889941304399,  0.3700540788845391,   0.3700540788845392)
    ,@("FAPs", "Angpt1", "Itga5", "sCs",  3, 1,                  12.95459633333333, 38.863789,           0.8806046431437031,  0.8806046431437031,  3, 1, 4.125957666666666, 12.377873,         0.06262628758663766,0.06262628758663766,53.45011606008855,  481.0510445407969,  0.05514899963164598,  0.05514899963164598)
    ,@("sCs",  "Angpt1", "Itga5", "ECs",  3, 1,                  1.590497666666667, 4.771493,            0.1081160380560855,  0.1081160380560855,  3, 1, 34.07074633333333, 102.212239,        0.5171464495142372, 0.5171464495142373, 54.18944254475854,  487.7049829028269,  0.05591182521625077,  0.05591182521625078)
    ,@("sCs",  "Angpt1", "Itga5", "FAPs", 3, 1,                  1.590497666666667, 4.771493,            0.1081160380560855,  0.1081160380560855,  3, 1, 27.685497,          83.056491,         0.420227262899125,  0.4202272628991251, 44.03371837900699,  396.3034654110629,  0.04543330674780646,  0.04543330674780646)
    ,@("sCs",  "Angpt1", "Itga5", "sCs",  3, 1,                  1.590497666666667, 4.771493,            0.1081160380560855,  0.1081160380560855,  3, 1, 4.125957666666666, 12.377873,         0.06262628758663766,0.06262628758663766,6.562326041598777,  59.06093437438899,  0.006770906092028272,  0.006770906092028272)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
